$wb = $excel.ActiveWorkbook
foreach ($ws in $wb.Worksheets) {
    foreach ($addr in @("D5","E5","F5","G5","D7","E7","F7","G7","D8","E8","F8","G8")) {
        $cell = $ws.Range($addr)
        $val = $cell.Value2
        if ($val -ne 0) {
            $cell.Value2 = $val / 1000000
        }
    }
}
